# TC32_Canine_Filter_Breed-IrishWlfhnd.xlsx - "Fixed ICDC breed all testcases"
#
# The StatQuery column (column C) on the "startup" sheet held a stale/broken
# Cypher query (missing program/study file counts, malformed `IN[...]` syntax).
# Replace it with the corrected query that returns Programs/Studies/Cases/
# Samples/Case Files/Study Files counts, for every data row (Cases, Samples,
# Files tabs).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Irish Wolfhound']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Cells.Item(2, 3).Value = $newStatQuery
$ws.Cells.Item(3, 3).Value = $newStatQuery
$ws.Cells.Item(4, 3).Value = $newStatQuery

# Match the author's final view state: zoomed to 100% (was 55%) and the
# active selection left on B4 (was B2), with the sheet scrolled back to A1.
$ws.Range("A1").Select()
$excel.ActiveWindow.Zoom = 100
$ws.Range("B4").Select()
